$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 3 with the prior check-in record now marked OUT
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = 303
$ws.Range("B3").Value = 6
$ws.Range("C3").Value = "王*佑"
$ws.Range("D3").Value = "2023-07-04 15:27:23"
$ws.Range("E3").Value = "2023-07-04 15:42:50"
$ws.Range("F3").Value = "OUT"

# Update row 2 (new check-in) with the refreshed timestamp
$ws.Range("D2").Value = "2023-07-04 16:32:08"
